$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.247.31'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.87%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.902.17'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.13%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.28%  '

# Row 5
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.74'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.86%  '

# Row 6
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.690'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +9.03%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.30%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.79'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.09%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.349'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.89%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.01'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +11.94%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0725'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.95%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0994'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.12%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.177.74'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.13%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.31'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.10%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.711'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.64%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.909.05'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.85%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.83'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.07%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '35.276.56'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.94%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.42'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.35%  '

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.73%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '240.96'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.17%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.59'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.40%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.83'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.40%  '

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.31%  '

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.54%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +12.63%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.50'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.11%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.55'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.61%  '

# Row 29
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.36'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.67%  '

# Row 30
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.131'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.04%  '

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.42%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.962'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.32%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0571'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.26%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.22%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.49%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.78'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.53%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.03'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.37%  '

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.50%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0684'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +15.20%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.10'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.24%  '

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.57%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.16'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +6.11%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '90.41'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.09%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.345.35'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.94%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.44'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.70%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '47.08'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.37%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.60'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.80%  '

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.01%  '

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.98%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.57'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.39%  '
